# aggiornamento fino a 1/09/2021
# Append new daily rows (358-366) to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row number, date serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$newRows = @(
    @(358, 44432, 1, 2, 24.30724355858046),
    @(359, 44433, 0, 2, 24.30724355858046),
    @(360, 44434, 0, 2, 24.30724355858046),
    @(361, 44435, 0, 1, 12.15362177929023),
    @(362, 44436, 0, 1, 12.15362177929023),
    @(363, 44437, 0, 1, 12.15362177929023),
    @(364, 44438, 0, 1, 12.15362177929023),
    @(365, 44439, 2, 2, 24.30724355858046),
    @(366, 44440, 0, 2, 24.30724355858046)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $dateSerial = $r[1]
    $nuoviPos = $r[2]
    $sommaMobile = $r[3]
    $somma100k = $r[4]

    # Column A carries the same date style (s="2") as the row above it -
    # copy formats from the last existing data row, then overwrite the value.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)
    $ws.Range("A" + $rowNum).Value = $dateSerial

    $ws.Range("B" + $rowNum).Value = $nuoviPos
    $ws.Range("C" + $rowNum).Value = $sommaMobile
    $ws.Range("D" + $rowNum).Value = $somma100k
}
